$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the DC cable length matrix (symmetric cluster-to-cluster distances, km)
# Previously nonzero cluster_0<->cluster_1, cluster_1<->cluster_2,
# cluster_2<->cluster_4, cluster_0<->cluster_7, cluster_3<->cluster_7 links
# are replaced by a new set of links:
# cluster_0<->cluster_2 (388), cluster_0<->cluster_6 (690),
# cluster_5<->cluster_6 (382)

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 388
$ws.Range("H2").Value = 690
$ws.Range("I2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 388
$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("I5").Value = 0

$ws.Range("D6").Value = 0
$ws.Range("G6").Value = 382

$ws.Range("F7").Value = 382

$ws.Range("B8").Value = 690

$ws.Range("B9").Value = 0
$ws.Range("E9").Value = 0
